# Apply updated crypto price / 1h-volume-change values to sheet1.
# Source: commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on Price (D) / Volume(1h) (E) cells before
# assigning values so strings that look numeric (e.g. "1.000", "0.9999")
# are stored verbatim instead of being auto-coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.844.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.800.44'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.96'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4677'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.54%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07373'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8687'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.24%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.835.46'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.351'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.35'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.481'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07025'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.840.69'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.292'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.60'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.008.07'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.890'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.59'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.30'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.140'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -8.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.256'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.00'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08935'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7566'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.65%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.147'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.451'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.0000'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.099'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01954'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.925'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.226'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.379'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5283'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1657'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.485'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4996'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.32'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.03'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9999'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.662'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06286'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.97%  '
